$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 453
$ws1.Range("F3").Value = 5463
$ws1.Range("F4").Value = 389
$ws1.Range("F6").Value = 79
$ws1.Range("F9").Value = 520

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 453
$ws4.Range("F3").Value = 5463
$ws4.Range("F4").Value = 389
$ws4.Range("F7").Value = 79
$ws4.Range("F11").Value = 520
